# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# Updates the "Estado de Cuenta" worker detail table (rows 16-21) on
# Hoja1: some prior debtor rows are removed/reordered and new ones
# (LUIS ENRIQUE JULIO LUQUEZ with periods 1712/1801/1802, replacing the
# CARLOS ENRIQUE VANEGAS CAÑATE single row which moves to the bottom)
# are incorporated, per the refreshed source database.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Columns: B = Tipo Doc Trabajador, C = N Doc Trabajador, D = Nombre Trabajador,
#          E = Periodo Mora, F = Valor Mora, G = Salario Basico
# (B stays "CC" for every data row; it is left untouched.)

$ws.Range("C16").Value = "1002322141"
$ws.Range("D16").Value = "YONER LUIS PEREZ TORRES"
$ws.Range("E16").Value = "1607"
$ws.Range("F16").Value = 1839
$ws.Range("G16").Value = 689455

$ws.Range("C17").Value = "1143333674"
$ws.Range("D17").Value = "JAINER RAFAEL TORRES JULIO"
$ws.Range("E17").Value = "1607"
$ws.Range("F17").Value = 1839
$ws.Range("G17").Value = 689455

$ws.Range("C18").Value = "73100085"
$ws.Range("D18").Value = "LUIS ENRIQUE JULIO LUQUEZ"
$ws.Range("E18").Value = "1802"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 781242

$ws.Range("C19").Value = "73100085"
$ws.Range("D19").Value = "LUIS ENRIQUE JULIO LUQUEZ"
$ws.Range("E19").Value = "1801"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 781242

$ws.Range("C20").Value = "73100085"
$ws.Range("D20").Value = "LUIS ENRIQUE JULIO LUQUEZ"
$ws.Range("E20").Value = "1712"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 781242

$ws.Range("C21").Value = "86677305"
$ws.Range("D21").Value = "CARLOS ENRIQUE VANEGAS CAÑATE"
$ws.Range("E21").Value = "1609"
$ws.Range("F21").Value = 3677
$ws.Range("G21").Value = 689455
